$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("dataText")

# Remove the trailing fullstop from every cell in column B (rows 2-24) whose
# text ends with "now on." - this fixes the double fullstop issue described
# in the commit message ("...North East from now on." -> "...North East from now on")
for ($r = 2; $r -le 24; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    $val = $cell.Value2
    if ($val -ne $null -and $val.EndsWith("now on.")) {
        $cell.Value2 = $val.Substring(0, $val.Length - 1)
    }
}

# Update the active selection on the sheet from B24 to B25
$ws.Range("B25").Select()
